# The commit inserts one new price-observation row (row 448) into the
# "Hortaliza, Terminal Hortofrutícola Agro Chillán - Tomate" data sheet,
# pushing the previously existing rows 448-530 down to 449-531.
#
# Result: dimension grows from A1:R530 to A1:R531, and the whole block of
# rows from 448 onward is shifted down by one row, with the freshly
# inserted row 448 holding a new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 448; Excel shifts rows 448:530 -> 449:531
# and preserves per-column formatting (date style on column D, etc.)
$ws.Rows("448:448").Insert()

# Populate the newly inserted row 448 with the new observation.
$ws.Range("A448").Value = 7
$ws.Range("B448").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C448").Value = "Ñuble"
$ws.Range("D448").Value = 44816
$ws.Range("E448").Value = 16
$ws.Range("F448").Value = 100112020
$ws.Range("G448").Value = "Tomate"
$ws.Range("H448").Value = "Larga vida"
$ws.Range("I448").Value = "Primera"
$ws.Range("J448").Value = 400
$ws.Range("K448").Value = 10000
$ws.Range("L448").Value = 11000
$ws.Range("M448").Value = 10500
$ws.Range("N448").Value = "$/bandeja 18 kilos"
$ws.Range("O448").Value = "Región de Arica y Parinacota"
$ws.Range("P448").Value = 583
$ws.Range("Q448").Value = 18
$ws.Range("R448").Value = "Hortaliza"
